$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38 (holomorphic): update the "C" and "D" notes ------------------
# Old:  C38 = "복소해석", D38 = "복소해석으로 통일함"
# New:  C38 = "복소해석, 해석적", D38 = "복소해석으로 통일함 형용사로 해석적 (서울대 교재 참고)"
$ws.Range("C38").Value() = "복소해석, 해석적"
$ws.Range("D38").Value() = "복소해석으로 통일함 형용사로 해석적 (서울대 교재 참고)"

# --- New rows for Chapter 4 (Laurent series) terms ------------------------
$ws.Range("A99").Value() = "removable singularity"
$ws.Range("B99").Value() = "제거가능한 특이점"

$ws.Range("A100").Value() = "essential singularity"
$ws.Range("B100").Value() = "본질적 특이점"
$ws.Range("D100").Value() = "용어사전"

$ws.Range("A101").Value() = "Casorati"
$ws.Range("B101").Value() = "카소라티"
$ws.Range("D101").Value() = "임의로 번역"

# --- Update the view's active cell / selection to match the scrolled state
$ws.Range("A102").Select()
